# Update generated figures ("想去人数" / "最低票价") across the four
# sheets of the workbook, matching the refreshed scrape output.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" -----------------------------------------------------
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 2891
$ws.Range("F3").Value = 21229
$ws.Range("G3").Value = "不可售"
$ws.Range("F5").Value = 3012
$ws.Range("F8").Value = 517
$ws.Range("F10").Value = 282
$ws.Range("F13").Value = 120
$ws.Range("F14").Value = 523
$ws.Range("F15").Value = 184
$ws.Range("F19").Value = 69
$ws.Range("F22").Value = 38
$ws.Range("F23").Value = 128

# --- Sheet "演出" -----------------------------------------------------
$ws = $wb.Worksheets.Item("演出")
$ws.Range("G2").Value = 138
$ws.Range("F3").Value = 34
$ws.Range("F12").Value = 102

# --- Sheet "本地生活" --------------------------------------------------
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 6134
$ws.Range("F3").Value = 702
$ws.Range("F4").Value = 700
$ws.Range("F5").Value = 1634
$ws.Range("F6").Value = 57

# --- Sheet "全部类型" (combined view mirrors the same rows) -----------
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 6134
$ws.Range("F3").Value = 702
$ws.Range("F4").Value = 700
$ws.Range("F5").Value = 1634
$ws.Range("F6").Value = 2891
$ws.Range("F7").Value = 21229
$ws.Range("G7").Value = "不可售"
$ws.Range("G8").Value = 138
$ws.Range("F9").Value = 34
$ws.Range("F13").Value = 3012
$ws.Range("F16").Value = 57
$ws.Range("F18").Value = 517
$ws.Range("F20").Value = 282
$ws.Range("F26").Value = 120
$ws.Range("F29").Value = 523
$ws.Range("F30").Value = 102
$ws.Range("F31").Value = 184
$ws.Range("F39").Value = 69
$ws.Range("F44").Value = 38
$ws.Range("F50").Value = 128
